$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1:C1").EntireColumn.Insert()
$col = $ws.Columns.Item(2)
$col.ColumnWidth = 35.855
$col.WrapText = $true
$col.Borders.Item(8).LineStyle = 1
Write-Host "done"
